$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# C40 currently holds the text "NA" -> clear its content but keep it as an
# (empty) text cell, like the other "Rien ne nous concerne..." rows whose
# "Numero de page" column is blank text rather than a totally empty cell.
# A leading apostrophe forces text-typed content; with nothing after it the
# cell becomes an empty string. Then drop the formatting Excel stamped on
# it while doing that, so no stray style sticks to the cell.
$ws.Range("C40").Value = "'"
$ws.Range("C40").ClearFormats()

# Add new row 41 with data.
# A41 looks like a date ("2025-04-04"); force it to stay text like the
# other Date-column cells (avoid Excel's automatic date inference), then
# drop the temporary formatting again so no stray style sticks to the cell.
$ws.Range("A41").NumberFormat = "@"
$ws.Range("A41").Value = "2025-04-04"
$ws.Range("A41").ClearFormats()

$ws.Range("B41").Value = "ruissellement"
$ws.Range("C41").Value = 56
$ws.Range("D41").Value = 1
